# LR3 - seems to be ready. Pattern composite
#
# The deck's Slide Master and every Slide Layout carry a fixed
# "Date placeholder" (Insert > Header & Footer > Date and time > Fixed)
# that currently reads "01-Apr-22". Bring it forward to "18-Apr-22"
# everywhere it appears (the master plus all custom layouts), same as
# re-applying the Header & Footer dialog with the updated date to all.

$p = $ppt.ActivePresentation
$newDate = "18-Apr-22"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster

# Slide Master's own date placeholder.
Update-DatePlaceholder $master.Shapes

# Every slide layout owned by the master.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
